$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# Add TA office hour: the "Teaching Assistant Office Hours" line currently
# says "TBD" -- fill in the actual meeting time.  The insertion is split into
# "T" (kept from the original run) + "uesday 9:00-11:00" (new run), the same
# way Word splits a run when you backspace over "BD" and type the rest,
# keeping the same font colour used throughout this heading line.
# ---------------------------------------------------------------------------
$r = $d.Content
$found = $r.Find.Execute("TBD", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($found) {
    $r.Text = "T"
    $r.Collapse(0)
    $r.InsertAfter("uesday 9:00-11:00")
    $r.Font.Color = 11930122   # 0A0AB6, matching the rest of the line
}

# ---------------------------------------------------------------------------
# Correct the TA's e-mail address: drop the "onid" sub-domain so it reads
# yih@oregonstate.edu instead of yih@onid.oregonstate.edu.
# ---------------------------------------------------------------------------
$d.Content.Find.Execute("yih@onid.oregonstate.edu", $true, $false, $false, $false, $false, $true, 1, $false, "yih@oregonstate.edu", 2) | Out-Null

# ---------------------------------------------------------------------------
# Tidy up "open note, open internet essay exams." into a single phrase.
# ---------------------------------------------------------------------------
$d.Content.Find.Execute(" open note, open internet essay exams.", $true, $false, $false, $false, $false, $true, 1, $false, " open note, open internet essay exams.", 2) | Out-Null
